$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contact info")

# Insert a new blank column before column E ("Email"), shifting the
# existing "Email" -> F and "Interesting comments" -> G.
$ws.Columns("E:E").Insert()

# Fill in the (now-empty) "Interesting comments" column with notes
# gathered from today's MSTP meetings.
$ws.Range("E4").Value = "very offputting. Almost rude"
$ws.Range("E28").Value = "super kind, invested in my future. Email about age"

# The two mailto hyperlinks that used to live in column E now live in
# column F (shifted by the column insert) - move them over.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F7"), "mailto:zachary.d.ewell@vanderbilt.edu")
$ws.Hyperlinks.Add($ws.Range("F12"), "mailto:schapman@health.ucsd.edu")

# Re-adding the hyperlinks nudges the cell style; restore the original
# blue-highlight "Hyperlink" look those two cells had before the move.
$ws.Range("F7").Font.Underline = $true
$ws.Range("F12").Font.Underline = $true

# Leave the cursor where it was last left after typing the second note.
[void]$ws.Range("E29").Select()
